$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "B4"  = 8.838199999999995
    "B6"  = 5.047000000000006
    "B7"  = 5.643000000000004
    "B16" = 4.935799999999999
    "B20" = 9.515299999999986
    "B28" = 6.035900000000002
    "B29" = 5.125200000000004
    "B32" = 6.614799999999999
    "B40" = 9.146699999999992
    "B46" = 5.646600000000006
    "B51" = 5.532499999999998
    "B52" = 5.377099999999997
    "B57" = 5.030399999999995
    "B59" = 4.749199999999999
    "B62" = 5.503999999999999
    "B66" = 5.847899999999995
    "B73" = 8.9176
    "B74" = 9.054699999999995
    "B92" = 4.768299999999998
    "B100" = 4.957600000000005
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
